# Contoso Chai Tea market trends 2023 - rename "Chai Tea"-specific header
# labels to generic "Tea" labels and bold the (white) header text, matching
# the table header row B1:F1 (A1 / "日期" was already bold).
#
# Updating Range.Value on a table header cell also renames the matching
# ListObject/ListColumn automatically, so editing the cells is sufficient
# to keep the table definition (Table1) in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "总茶销售量（单位）"
$ws.Range("C1").Value = "手工制茶销售量（单位）"
$ws.Range("D1").Value = "预制茶销售量（单位）"
$ws.Range("F1").Value = "线上茶搜索量"

# Bold every header cell's text (A1 already was bold). Re-assert the white
# font color too, since replacing .Value drops the prior explicit run
# formatting and would otherwise fall back to the default (black) font.
$header = $ws.Range("B1:F1")
$header.Font.Bold = $true
$header.Font.Color = 16777215
